# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, and swaps the EnergySwap / Filecoin row order (rows 44-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.379.55"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "3.508.42"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("D5").Value = "'591.13"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "'7.64"
$ws.Range("E9").Value = "  +6.20%  "

$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").Value = "'0.388"
$ws.Range("E11").Value = "  +3.36%  "

$ws.Range("D12").Value = "4.106.73"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").Value = "3.508.36"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "64.364.49"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "'25.72"
$ws.Range("E17").Value = "  +2.81%  "

$ws.Range("D18").Value = "'10.06"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  +2.44%  "

$ws.Range("D20").Value = "'13.62"
$ws.Range("E20").Value = "  -1.24%  "

$ws.Range("D21").Value = "'395.08"
$ws.Range("E21").Value = "  +2.71%  "

$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").Value = "3.647.86"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "'74.71"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").Value = "'5.73"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  +2.73%  "

$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  -1.59%  "

$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("D31").Value = "'8.26"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("E32").Value = "  -6.37%  "

$ws.Range("E33").Value = "  +6.15%  "

$ws.Range("D34").Value = "3.537.63"
$ws.Range("E34").Value = "  +0.44%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'23.37"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("E37").Value = "  +0.91%  "

$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("E39").Value = "  +0.80%  "

$ws.Range("D40").Value = "'167.24"
$ws.Range("E40").Value = "  +2.30%  "

$ws.Range("D41").Value = "'0.0787"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.45"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'25.00"
$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("D46").Value = "'1.66"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("E47").Value = "  -3.25%  "

$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").Value = "2.379.51"
$ws.Range("E49").Value = "  -3.85%  "

$ws.Range("D50").Value = "'0.896"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").Value = "'0.0259"
$ws.Range("E51").Value = "  -0.30%  "
